# Update the "Förändrad" (Changed) date column (C) for rows 2-7
# from 2023-09-01 (45170) to 2023-09-05 (45174), matching the
# automatic update reflected in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = 45174
